$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "state" (B) and "city" (C) columns' content - they were removed entirely.
$ws.Range("B1:C3").ClearContents()

# Update A3 value from 30000 to 16000
$ws.Range("A3").Value = 16000

# Update the active selection to match the target (E9)
$ws.Range("E9").Select()
